$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update data for rows 2-37 (columns A-E) ---
$data = @(
    ,@("Borderlands 3", 0.85, 99074, 8.99, 0.85)
    ,@("BioShock Infinite", 0.93, 100694, 7.49, 0.75)
    ,@("Hell Let Loose", 0.84, 73081, 29.24, 0.35)
    ,@("Trepang2", 0.9399999999999999, 5445, 20.09, 0.33)
    ,@("Warhammer 40,000: Boltgun", 0.91, 9774, 14.95, 0.32)
    ,@("COCOON", 0.95, 3262, 16.09, 0.3)
    ,@("Roboquest", 0.95, 9633, 18.74, 0.25)
    ,@("Fireworks Mania - An Explosive Simulator", 0.96, 4813, 7.99, 0.2)
    ,@("Tom Clancy's Rainbow Six® Siege", 0.86, 1032092, 19.99, 0)
    ,@("Team Fortress 2", 0.93, 1022451, 0, 0)
    ,@("BattleBit Remastered", 0.9, 107017, 14.79, 0)
    ,@("Deep Rock Galactic", 0.97, 215574, 29.99, 0)
    ,@("ULTRAKILL", 0.98, 83731, 24.5, 0)
    ,@("Inscryption", 0.96, 96881, 19.99, 0)
    ,@("Gunfire Reborn", 0.93, 79152, 16.79, 0)
    ,@("A Little to the Left", 0.91, 6209, 14.99, 0)
    ,@("Metal: Hellsinger", 0.96, 10895, 29.99, 0)
    ,@("Left 4 Dead 2", 0.97, 614416, 9.75, 0)
    ,@("The Talos Principle 2", 0.95, 6562, 28.99, 0)
    ,@("Escape Simulator", 0.9399999999999999, 11499, 14.99, 0)
    ,@("Portal 2", 0.98, 310744, 9.75, 0)
    ,@("We Were Here Forever", 0.91, 10575, 17.99, 0)
    ,@("Chants of Sennaar", 0.98, 8513, 19.99, 0)
    ,@("GROUND BRANCH", 0.91, 15076, 24.99, 0)
    ,@("PAYDAY 2", 0.89, 423660, 9.99, 0)
    ,@("Far Cry® 4", 0.83, 46590, 29.99, 0)
    ,@("Jusant", 0.96, 1347, 24.99, 0)
    ,@("SCP: 5K", 0.82, 6685, 19.5, 0)
    ,@("Metro Exodus", 0.89, 86868, 29.99, 0)
    ,@("We Were Here Together", 0.85, 10331, 12.99, 0)
    ,@("MechWarrior 5: Mercenaries", 0.84, 8316, 29.99, 0)
    ,@("Call of Duty: World at War", 0.92, 39879, 19.99, 0)
    ,@("Dorfromantik", 0.96, 23024, 12.99, 0)
    ,@("Squad 44", 0.8, 15920, 28, 0)
    ,@("TUNIC", 0.92, 9826, 28.99, 0)
    ,@("Starship Troopers: Extermination", 0.89, 23242, 28.99, 0)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowVals = $data[$i]
    $ws.Range("A$row").Value2 = $rowVals[0]
    $ws.Range("B$row").Value2 = $rowVals[1]
    $ws.Range("C$row").Value2 = $rowVals[2]
    $ws.Range("D$row").Value2 = $rowVals[3]
    $ws.Range("E$row").Value2 = $rowVals[4]
}

# --- Remove rows 38-50 and their hyperlinks, keep hyperlinks for rows 2-37 ---
$keepUrls = @(
    ,"https://store.steampowered.com/app/673610/Airport_CEO/"
    ,"https://store.steampowered.com/app/2287220/F1_Manager_2023/"
    ,"https://store.steampowered.com/app/1080020/Electrician_Simulator/"
    ,"https://store.steampowered.com/app/1480560/Lawn_Mowing_Simulator/"
    ,"https://store.steampowered.com/app/1714250/Pets_Hotel/"
    ,"https://store.steampowered.com/app/1137750/Farmers_Life/"
    ,"https://store.steampowered.com/app/1214470/Hotel_Renovator/"
    ,"https://store.steampowered.com/app/1899350/Trade_Bots_A_Technical_Analysis_Simulation/"
    ,"https://store.steampowered.com/app/849100/Alaskan_Road_Truckers/"
    ,"https://store.steampowered.com/app/2000160/Music_Store_Simulator/"
    ,"https://store.steampowered.com/app/2529170/Storage_Hustle/"
    ,"https://store.steampowered.com/app/1139980/Travellers_Rest/"
    ,"https://store.steampowered.com/app/362620/Software_Inc/"
    ,"https://store.steampowered.com/app/2460920/ACRES/"
    ,"https://store.steampowered.com/app/1244910/My_Supermarket/"
    ,"https://store.steampowered.com/app/1059900/Tribe_Primitive_Builder/"
    ,"https://store.steampowered.com/app/2479290/Computer_Repair_Shop/"
    ,"https://store.steampowered.com/app/1270580/Mind_Over_Magic/"
    ,"https://store.steampowered.com/app/997010/Police_Simulator_Patrol_Officers/"
    ,"https://store.steampowered.com/app/986130/Shadows_of_Doubt/"
    ,"https://store.steampowered.com/app/573090/Stormworks_Build_and_Rescue/"
    ,"https://store.steampowered.com/app/382310/Eco/"
    ,"https://store.steampowered.com/app/2653790/The_Exit_8/"
    ,"https://store.steampowered.com/app/371970/Barony/"
    ,"https://store.steampowered.com/app/777390/Flyout/"
    ,"https://store.steampowered.com/app/1388770/Cruelty_Squad/"
    ,"https://store.steampowered.com/app/2248760/Car_For_Sale_Simulator_2023/"
    ,"https://store.steampowered.com/app/2559270/Gym_Simulator_24/"
    ,"https://store.steampowered.com/app/621060/PC_Building_Simulator/"
    ,"https://store.steampowered.com/app/2546690/Tram_Simulator_Urban_Transit/"
    ,"https://store.steampowered.com/app/1705180/Gunner_HEAT_PC/"
    ,"https://store.steampowered.com/app/24780/SimCity_4_Deluxe_Edition/"
    ,"https://store.steampowered.com/app/1999360/Placid_Plastic_Duck_Simulator/"
    ,"https://store.steampowered.com/app/1150760/Gloomwood/"
    ,"https://store.steampowered.com/app/2477090/Mosa_Lina/"
    ,"https://store.steampowered.com/app/1122340/Chef_Life_A_Restaurant_Simulator/"
)

$ws.Hyperlinks.Delete()
$ws.Rows("38:50").Delete()

for ($i = 0; $i -lt $keepUrls.Length; $i++) {
    $row = $i + 2
    $ws.Hyperlinks.Add($ws.Range("A$row"), $keepUrls[$i]) | Out-Null
}

Write-Host "Final UsedRange: $($ws.UsedRange.Address())"
Write-Host "Final Hyperlinks count: $($ws.Hyperlinks.Count)"
